$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 119, shifting rows 119-181 down to 120-182.
# This automatically carries the old row 181 data into the new row 182,
# matching the target diff (row 182 == old row 119..181's trailing record).
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the latest weekly record.
$ws.Range("A119").Value = 8
$ws.Range("B119").Value = "Terminal La Palmera de La Serena"
$ws.Range("C119").Value = "Coquimbo"
$ws.Range("D119").Value = 44518
$ws.Range("E119").Value = 4
$ws.Range("F119").Value = 100112012
$ws.Range("G119").Value = "Espinaca"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 3100
$ws.Range("K119").Value = 400
$ws.Range("L119").Value = 500
$ws.Range("M119").Value = 450
$ws.Range("N119").Value = "$/atado 300 a 500 gramos"
$ws.Range("O119").Value = "Provincia del Elquí"
$ws.Range("P119").Value = 900
$ws.Range("Q119").Value = 0.5
$ws.Range("R119").Value = "Hortaliza"
